# Generate Report for Archive
#
# 1) The shared status string "Ready for handoff" becomes "In Translation".
#    That string is referenced from four cells across the workbook
#    (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2), so all four must be
#    updated so every reference to the old text is replaced.
# 2) The "zh-cn"/"de-de" status columns (and the mirrored columns on the
#    Overview sheet) are narrowed from ~17.22 to ~13.41 character-width
#    units.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Update the status text everywhere it appears -----------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- 2) Narrow the status columns -------------------------------------
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth  # column F

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth      # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth      # column C (Status)
